$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing hyperlink info (row + target URL) before shifting columns,
# since the engine does not relocate Hyperlink ranges automatically on
# column insert.
$hlInfo = @()
foreach ($hl in $ws.Hyperlinks) {
    $r = $hl.Range.Row
    $target = $hl.Address
    $hlInfo += [PSCustomObject]@{ Row = $r; Target = $target }
}

# Insert a new column before column E (Link), shifting the Link column (and
# its contents/styles) to column F.
$ws.Columns.Item(5).Insert()

# New column E header + values: "Category"
$ws.Range("E1").Value = "Category"

# Fill Category values for rows 2-101 based on column A (Picture Series number):
# 1 for series 0-14, 2 for series 15-19
for ($r = 2; $r -le 101; $r++) {
    $seriesVal = [double]$ws.Cells.Item($r, 1).Value2
    if ($seriesVal -le 14) {
        $ws.Cells.Item($r, 5).Value = 1
    } else {
        $ws.Cells.Item($r, 5).Value = 2
    }
}

# Rebuild hyperlinks so they point at column F (where the Link text now lives).
$ws.Hyperlinks.Delete()
foreach ($item in $hlInfo) {
    $target = $ws.Cells.Item($item.Row, 6)
    $ws.Hyperlinks.Add($target, $item.Target) | Out-Null
    # Hyperlinks.Add() re-stamps the style, which can create a duplicate (but
    # equivalent) style record; re-applying the named style collapses it back
    # onto the original shared "Hyperlink" style entry.
    $target.Style = "Hyperlink"
}

# Update selection to match target state
$ws.Range("E105").Select()
